$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44749
$ws.Range("M2").Value = 120

# Row 3
$ws.Range("D3").Value = 44405
$ws.Range("M3").Value = 50
$ws.Range("N3").Value = 1200
$ws.Range("O3").Value = 1200
$ws.Range("P3").Value = 1200
$ws.Range("S3").Value = 1200

# Row 4
$ws.Range("D4").Value = 44418
$ws.Range("M4").Value = 40
$ws.Range("N4").Value = 1200
$ws.Range("O4").Value = 1200
$ws.Range("P4").Value = 1200
$ws.Range("S4").Value = 1200

# Row 5
$ws.Range("D5").Value = 44811
$ws.Range("M5").Value = 60
$ws.Range("N5").Value = 2500
$ws.Range("O5").Value = 2500
$ws.Range("P5").Value = 2500
$ws.Range("S5").Value = 2500

# Row 6
$ws.Range("D6").Value = 44432
$ws.Range("M6").Value = 30
$ws.Range("N6").Value = 1300
$ws.Range("O6").Value = 1300
$ws.Range("P6").Value = 1300
$ws.Range("S6").Value = 1300

# Row 7
$ws.Range("D7").Value = 44812
$ws.Range("M7").Value = 50
$ws.Range("N7").Value = 2500
$ws.Range("O7").Value = 2500
$ws.Range("P7").Value = 2500
$ws.Range("S7").Value = 2500

# Row 8
$ws.Range("D8").Value = 44357
$ws.Range("M8").Value = 35
$ws.Range("N8").Value = 1000
$ws.Range("O8").Value = 1000
$ws.Range("P8").Value = 1000
$ws.Range("S8").Value = 1000

# Row 9
$ws.Range("D9").Value = 44424
$ws.Range("M9").Value = 50

# Row 10
$ws.Range("D10").Value = 44435
$ws.Range("M10").Value = 130

# Row 11
$ws.Range("D11").Value = 44763
$ws.Range("M11").Value = 50
$ws.Range("N11").Value = 2300
$ws.Range("O11").Value = 2300
$ws.Range("P11").Value = 2300
$ws.Range("S11").Value = 2300

# Row 12
$ws.Range("D12").Value = 44830
$ws.Range("M12").Value = 50

# Row 13
$ws.Range("D13").Value = 44431
$ws.Range("M13").Value = 100
$ws.Range("N13").Value = 1300
$ws.Range("O13").Value = 1300
$ws.Range("P13").Value = 1300
$ws.Range("S13").Value = 1300

# Row 14
$ws.Range("D14").Value = 44417
$ws.Range("M14").Value = 80

# Row 15
$ws.Range("D15").Value = 44760
$ws.Range("M15").Value = 80
$ws.Range("N15").Value = 2300
$ws.Range("O15").Value = 2300
$ws.Range("P15").Value = 2300
$ws.Range("S15").Value = 2300

# Row 16
$ws.Range("D16").Value = 44748
$ws.Range("M16").Value = 300

# Row 17
$ws.Range("D17").Value = 44473
$ws.Range("M17").Value = 120

# Row 18
$ws.Range("D18").Value = 44762
$ws.Range("N18").Value = 2300
$ws.Range("O18").Value = 2300
$ws.Range("P18").Value = 2300
$ws.Range("S18").Value = 2300

# Row 19
$ws.Range("D19").Value = 44753
$ws.Range("M19").Value = 160
$ws.Range("N19").Value = 2300
$ws.Range("O19").Value = 2300
$ws.Range("P19").Value = 2300
$ws.Range("S19").Value = 2300

# Row 20
$ws.Range("D20").Value = 44438
$ws.Range("M20").Value = 60

# Row 21
$ws.Range("D21").Value = 44343
$ws.Range("M21").Value = 60
$ws.Range("N21").Value = 1300
$ws.Range("O21").Value = 1300
$ws.Range("P21").Value = 1300
$ws.Range("S21").Value = 1300

# Row 22
$ws.Range("D22").Value = 44476
$ws.Range("M22").Value = 80
$ws.Range("N22").Value = 1200
$ws.Range("O22").Value = 1200
$ws.Range("P22").Value = 1200
$ws.Range("S22").Value = 1200
